# "Created Test Cases for Test Plans CRUD"
#
# The checklist's "In Qase" status option is retired - the four Test Plans
# CRUD rows (D2:D5), which still showed "In Qase", now show "Automated"
# instead (this also drops "In Qase" from the shared-string table since it
# becomes unused, shifting every later shared-string index down by one - the
# other rows' text is unaffected by that reindexing).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2:D5").Value = "Automated"

# B15 ("Create 3 Test Cases in 1 Suite and add 2 to Test Plan") carried a
# stray "apply bold font" flag even though it rendered with the regular
# (non-bold) font; clearing Bold normalizes it to the same format already
# used by the other centered/wrapped cells in column B (e.g. B2).
$ws.Range("B15").Font.Bold = $false

# The active selection moves from D16 to D6.
$ws.Range("D6").Select()
